# Added results for partial model
# Updates the sliding-window results (IPC PO predictions in column C,
# DELTA in column D, DELTA^2 in column E, plus the TOTAL/MSE summary
# rows) with the values produced by the partial model run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 29.15000000000001
$ws.Range("C2").Value = 29.39328575134277
$ws.Range("D2").Value = 0.2432857513427678
$ws.Range("E2").Value = 0.05918795680641502
$ws.Range("B3").Value = 29.34999999999999
$ws.Range("C3").Value = 29.74392509460449
$ws.Range("D3").Value = 0.3939250946044979
$ws.Range("E3").Value = 0.1551769801591626
$ws.Range("C4").Value = 29.60677146911621
$ws.Range("D4").Value = 0.2367714691162064
$ws.Range("E4").Value = 0.05606072858744668
$ws.Range("B5").Value = 29.53999999999999
$ws.Range("C5").Value = 29.0991268157959
$ws.Range("D5").Value = -0.4408731842040936
$ws.Range("E5").Value = 0.1943691645502567
$ws.Range("C6").Value = 29.46976280212402
$ws.Range("D6").Value = -0.08023719787597372
$ws.Range("E6").Value = 0.006438007922988162
$ws.Range("C7").Value = 29.63061332702637
$ws.Range("D7").Value = -0.1193866729736328
$ws.Range("E7").Value = 0.01425317768371315
$ws.Range("C8").Value = 30.01024436950684
$ws.Range("D8").Value = 0.1702443695068325
$ws.Range("E8").Value = 0.02898314534877893
$ws.Range("C9").Value = 29.95468711853027
$ws.Range("D9").Value = 0.1446871185302712
$ws.Range("E9").Value = 0.02093436226859274
$ws.Range("C10").Value = 29.93261528015137
$ws.Range("D10").Value = 0.01261528015136548
$ws.Range("E10").Value = 0.0001591452932974359
$ws.Range("C11").Value = 29.89081382751465
$ws.Range("D11").Value = -0.08918617248535554
$ws.Range("E11").Value = 0.00795417336258759
$ws.Range("B12").Value = 30.03999999999999
$ws.Range("C12").Value = 30.07125663757324
$ws.Range("D12").Value = 0.03125663757325015
$ws.Range("E12").Value = 0.0009769773923855127
$ws.Range("B13").Value = 30.21000000000001
$ws.Range("C13").Value = 30.06707572937012
$ws.Range("D13").Value = -0.1429242706298908
$ws.Range("E13").Value = 0.02042734713508626
$ws.Range("C14").Value = 30.21822166442871
$ws.Range("D14").Value = -0.001778335571287926
$ws.Range("E14").Value = 0.000003162477404107953
$ws.Range("C15").Value = 30.26735496520996
$ws.Range("D15").Value = -0.1126450347900345
$ws.Range("E15").Value = 0.01268890386284809
$ws.Range("C16").Value = 30.5389461517334
$ws.Range("D16").Value = 0.09894615173340071
$ws.Range("E16").Value = 0.009790340942849156
$ws.Range("C17").Value = 30.38221549987793
$ws.Range("D17").Value = -0.09778450012207429
$ws.Range("E17").Value = 0.009561808464123947
$ws.Range("C18").Value = 30.40595817565918
$ws.Range("D18").Value = -0.284041824340818
$ws.Range("E18").Value = 0.08067975797486013
$ws.Range("C19").Value = 30.4942684173584
$ws.Range("D19").Value = -0.2557315826416016
$ws.Range("E19").Value = 0.06539864236037829
$ws.Range("C20").Value = 30.59984397888184
$ws.Range("D20").Value = -0.3401560211181618
$ws.Range("E20").Value = 0.1157061187029393
$ws.Range("C21").Value = 30.71325492858887
$ws.Range("D21").Value = -0.2367450714111357
$ws.Range("E21").Value = 0.05604822883746372
$ws.Range("C22").Value = 31.08200645446777
$ws.Range("D22").Value = 0.06200645446777742
$ws.Range("E22").Value = 0.003844800395664554
$ws.Range("C23").Value = 31.26099967956543
$ws.Range("D23").Value = 0.1409996795654251
$ws.Range("E23").Value = 0.01988090963755257
$ws.Range("C24").Value = 31.31937217712402
$ws.Range("D24").Value = 0.0393721771240223
$ws.Range("E24").Value = 0.001550168331485385
$ws.Range("C25").Value = 31.23164939880371
$ws.Range("D25").Value = -0.1483506011962845
$ws.Range("E25").Value = 0.02200790087529905
$ws.Range("C26").Value = 31.44390678405762
$ws.Range("D26").Value = -0.1360932159423811
$ws.Range("E26").Value = 0.01852136342553957
$ws.Range("B27").Value = 31.65000000000001
$ws.Range("C27").Value = 31.86181449890137
$ws.Range("D27").Value = 0.2118144989013615
$ws.Range("E27").Value = 0.04486538194483487
$ws.Range("C28").Value = 32.47047805786133
$ws.Range("D28").Value = 0.5904780578613327
$ws.Range("E28").Value = 0.3486643368156913
$ws.Range("C29").Value = 32.40864944458008
$ws.Range("D29").Value = 0.128649444580077
$ws.Range("E29").Value = 0.0165506795907623
$ws.Range("C30").Value = 32.52600479125977
$ws.Range("D30").Value = 0.07600479125976278
$ws.Range("E30").Value = 0.005776728294440113
$ws.Range("B31").Value = 32.84999999999999
$ws.Range("C31").Value = 32.81044387817383
$ws.Range("D31").Value = -0.03955612182616619
$ws.Range("E31").Value = 0.001564686773926501
$ws.Range("B32").Value = 32.90000000000001
$ws.Range("C32").Value = 32.9732551574707
$ws.Range("D32").Value = 0.07325515747069744
$ws.Range("E32").Value = 0.005366318096056679
$ws.Range("B33").Value = 33.09999999999999
$ws.Range("C33").Value = 32.94234848022461
$ws.Range("D33").Value = -0.1576515197753849
$ws.Range("E33").Value = 0.02485400168748859
$ws.Range("B34").Value = 33.40000000000001
$ws.Range("C34").Value = 33.70544815063477
$ws.Range("D34").Value = 0.3054481506347599
$ws.Range("E34").Value = 0.093298572726195
$ws.Range("C35").Value = 33.68313980102539
$ws.Range("D35").Value = -0.01686019897461222
$ws.Range("E35").Value = 0.0002842663094635149
$ws.Range("B36").Value = 34.09999999999999
$ws.Range("C36").Value = 33.84911727905273
$ws.Range("D36").Value = -0.2508827209472599
$ws.Range("E36").Value = 0.0629421396699007
$ws.Range("B37").Value = 34.40000000000001
$ws.Range("C37").Value = 34.4184684753418
$ws.Range("D37").Value = 0.01846847534179119
$ws.Range("E37").Value = 0.0003410845814503493
$ws.Range("B38").Value = 34.90000000000001
$ws.Range("C38").Value = 34.99603652954102
$ws.Range("D38").Value = 0.09603652954100994
$ws.Range("E38").Value = 0.009223015006281275
$ws.Range("C39").Value = 35.76477432250977
$ws.Range("D39").Value = 0.4647743225097685
$ws.Range("E39").Value = 0.2160151708644143
$ws.Range("C40").Value = 36.0239372253418
$ws.Range("D40").Value = 0.323937225341794
$ws.Range("E40").Value = 0.1049353259621403
$ws.Range("C41").Value = 36.0843391418457
$ws.Range("D41").Value = -0.215660858154294
$ws.Range("E41").Value = 0.04650960573984653
$ws.Range("C42").Value = 36.63283157348633
$ws.Range("D42").Value = -0.167168426513669
$ws.Range("E42").Value = 0.02794528282305596
$ws.Range("C43").Value = 37.12261581420898
$ws.Range("D43").Value = -0.1773841857910128
$ws.Range("E43").Value = 0.03146514936874054
$ws.Range("B44").Value = 37.90000000000001
$ws.Range("C44").Value = 37.91189193725586
$ws.Range("D44").Value = 0.01189193725585369
$ws.Range("E44").Value = 0.000141418171697161
$ws.Range("C45").Value = 38.39858627319336
$ws.Range("D45").Value = -0.1014137268066406
$ws.Range("E45").Value = 0.01028474398481194
$ws.Range("B46").Value = 38.90000000000001
$ws.Range("C46").Value = 38.96462631225586
$ws.Range("D46").Value = 0.06462631225585369
$ws.Range("E46").Value = 0.004176560235791105
$ws.Range("B47").Value = 39.40000000000001
$ws.Range("C47").Value = 39.49670791625977
$ws.Range("D47").Value = 0.09670791625975994
$ws.Range("E47").Value = 0.009352421067304741
$ws.Range("B48").Value = 39.90000000000001
$ws.Range("C48").Value = 39.54153060913086
$ws.Range("D48").Value = -0.3584693908691463
$ws.Range("E48").Value = 0.1285003041900968
$ws.Range("B49").Value = 40.09999999999999
$ws.Range("C49").Value = 40.08893966674805
$ws.Range("D49").Value = -0.01106033325194744
$ws.Range("E49").Value = 0.0001223309716441342
$ws.Range("B50").Value = 40.59999999999999
$ws.Range("C50").Value = 40.62304306030273
$ws.Range("D50").Value = 0.02304306030274006
$ws.Range("E50").Value = 0.0005309826281157148
$ws.Range("B51").Value = 40.90000000000001
$ws.Range("C51").Value = 40.96342086791992
$ws.Range("D51").Value = 0.06342086791991619
$ws.Range("E51").Value = 0.004022206487715454

# TOTAL row
$ws.Range("C52").Value = 0.1406257629394361
$ws.Range("E52").Value = 2.178335986790985

# MSE row
$ws.Range("E53").Value = 0.04356671973581969
